$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.609.36'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '1.581.99'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.59'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.82%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '24.12'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('E10').Value = '  -1.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0592'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0893'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('D13').Value = '1.805.66'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('D14').Value = '1.583.06'
$ws.Range('E14').Value = '  -0.48%  '
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('D16').Value = '28.621.57'
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('E17').Value = '  -1.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '62.22'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '231.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.41'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('D21').Value = '0.0₃0692'
$ws.Range('E21').Value = '  -2.03%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('E23').Value = '  -3.56%  '
$ws.Range('E24').Value = '  -1.60%  '
$ws.Range('E25').Value = '  +5.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.11'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.05'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('E28').Value = '  -1.71%  '
$ws.Range('E29').Value = '  -2.31%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('E31').Value = '  +2.50%  '
$ws.Range('E32').Value = '  -1.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.21'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.14%  '
$ws.Range('E34').Value = '  -1.71%  '
$ws.Range('D35').Value = '1.395.90'
$ws.Range('E35').Value = '  -0.32%  '
$ws.Range('E36').Value = '  +4.10%  '
$ws.Range('E37').Value = '  -3.74%  '
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.65'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.80%  '
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.523'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.61%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.793'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.17%  '
$ws.Range('E44').Value = '  +1.69%  '
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('E46').Value = '  -2.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.961'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '63.20'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').Value = '1.719.27'
$ws.Range('E49').Value = '  -0.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.58'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.70%  '
$ws.Range('E51').Value = '  -2.50%  '
